$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 683, shifting existing rows 683-716 down to 684-717
$ws.Rows.Item(683).Insert()

# Populate the newly inserted row 683 with the new record
$ws.Cells.Item(683, 1).Value = 4
$ws.Cells.Item(683, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(683, 3).Value = "Los Lagos"
$ws.Cells.Item(683, 4).Value = 45267
$ws.Cells.Item(683, 4).NumberFormat = $ws.Cells.Item(684, 4).NumberFormat
$ws.Cells.Item(683, 5).Value = 10
$ws.Cells.Item(683, 6).Value = "Fruta"
$ws.Cells.Item(683, 7).Value = 100102
$ws.Cells.Item(683, 8).Value = "Cítricos"
$ws.Cells.Item(683, 9).Value = 100102006
$ws.Cells.Item(683, 10).Value = "Pomelo"
$ws.Cells.Item(683, 11).Value = "Start Ruby"
$ws.Cells.Item(683, 12).Value = "Primera"
$ws.Cells.Item(683, 13).Value = 50
$ws.Cells.Item(683, 14).Value = 14000
$ws.Cells.Item(683, 15).Value = 14000
$ws.Cells.Item(683, 16).Value = 14000
$ws.Cells.Item(683, 17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(683, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(683, 19).Value = 1000
$ws.Cells.Item(683, 20).Value = 14
